# Apply the changes described in the commit:
#  * Add global ligreto parameters (column widths adjusted to a uniform size)
#  * Add "missing" placeholder strings for rows where the right-hand side of
#    the stream-join had no matching record, and shrink the _FilterDatabase
#    defined name down to just the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
# Column B becomes a "default" sized column (stored width 8.0) while
# columns C:H all become a uniform, wider column (stored width ~17.59).
$ws.Columns("B").ColumnWidth = 7.09
$ws.Columns("C:H").ColumnWidth = 16.83

# --- Fill in "<<missing>>" markers for unmatched join rows --------------
# Row 5 : right-hand side (F:H) missing
$ws.Range("E5").Copy()
$ws.Range("F5:H5").PasteSpecial(-4122)
$ws.Range("F5:H5").Value = "<<missing>>"

# Row 6 : left-hand side (C:E) missing
$ws.Range("F6").Copy()
$ws.Range("C6:E6").PasteSpecial(-4122)
$ws.Range("C6:E6").Value = "<<missing>>"

# Row 7 : right-hand side (F:H) missing
$ws.Range("E7").Copy()
$ws.Range("F7:H7").PasteSpecial(-4122)
$ws.Range("F7:H7").Value = "<<missing>>"

# Row 8 : left-hand side (C:E) missing
$ws.Range("F8").Copy()
$ws.Range("C8:E8").PasteSpecial(-4122)
$ws.Range("C8:E8").Value = "<<missing>>"

$excel.CutCopyMode = 0

# --- Shrink the hidden _FilterDatabase defined name ----------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "full!_FilterDatabase") {
        $n.RefersTo = "=full!`$B`$2:`$H`$2"
    }
}
